$p = $ppt.ActivePresentation

# --- Update the "datetimeFigureOut" date placeholders from 01.11.2012 to 06.11.2012 ---

# 1) Slide Master "Date Placeholder 3" (ppt/slideMasters/slideMaster1.xml)
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "01.11.2012") {
        $shp.TextFrame.TextRange.Text = "06.11.2012"
    }
}

# 2) Slide Layout used by the title slide, "Date Placeholder 3" (ppt/slideLayouts/slideLayout13.xml)
$titleSlide = $p.Slides.Item(1)
$layout = $titleSlide.CustomLayout
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $shp = $layout.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "01.11.2012") {
        $shp.TextFrame.TextRange.Text = "06.11.2012"
    }
}

# 3) Notes Master "Date Placeholder 2" (ppt/notesMasters/notesMaster1.xml)
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "01.11.2012") {
        $shp.TextFrame.TextRange.Text = "06.11.2012"
    }
}

# --- Remove the vk.com link textbox from the front (title) slide ---
for ($i = $titleSlide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $titleSlide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "*vk.com*") {
        $shp.Delete()
    }
}
